$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("G6").Value = 2.1
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 3.8
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 1.91
$ws.Range("L6").Value = 4.5
$ws.Range("M6").Value = 1.11
$ws.Range("N6").Value = 6.5
$ws.Range("Q6").Value = 2.6
$ws.Range("R6").Value = 1.48
$ws.Range("S6").Value = 5.5
$ws.Range("T6").Value = 1.14
$ws.Range("W6").Value = 2.1
$ws.Range("X6").Value = 1.67
$ws.Range("Z6").Value = 9
$ws.Range("AA6").Value = 10
$ws.Range("AB6").Value = 19
$ws.Range("AC6").Value = 21
$ws.Range("AF6").Value = 6
$ws.Range("AG6").Value = 19
$ws.Range("AI6").Value = 8
$ws.Range("AJ6").Value = 17
$ws.Range("AL6").Value = 41
$ws.Range("AP6").Value = 2
$ws.Range("AQ6").Value = 1.85
$ws.Range("AR6").Value = 4.3
$ws.Range("AS6").Value = 1.22

# Row 7
$ws.Range("G7").Value = 2.5
$ws.Range("J7").Value = 3.5
$ws.Range("L7").Value = 4.33
$ws.Range("M7").Value = 1.17
$ws.Range("N7").Value = 5
$ws.Range("U7").Value = 1.73
$ws.Range("W7").Value = 2.5
$ws.Range("X7").Value = 1.5
$ws.Range("Z7").Value = 10
$ws.Range("AA7").Value = 12
$ws.Range("AB7").Value = 26
$ws.Range("AC7").Value = 29
$ws.Range("AD7").Value = 51
$ws.Range("AE7").Value = 4.75

# Row 15
$ws.Range("H15").Value = 5.25
$ws.Range("I15").Value = 9.5
$ws.Range("J15").Value = 1.77
$ws.Range("K15").Value = 2.4
$ws.Range("M15").Value = 1.04
$ws.Range("N15").Value = 13
$ws.Range("O15").Value = 1.25
$ws.Range("P15").Value = 3.75
$ws.Range("Q15").Value = 1.82
$ws.Range("R15").Value = 1.92
$ws.Range("S15").Value = 3.25
$ws.Range("T15").Value = 1.33
$ws.Range("U15").Value = 1.36
$ws.Range("V15").Value = 3
$ws.Range("W15").Value = 2.5
$ws.Range("X15").Value = 1.5
$ws.Range("Y15").Value = 6
$ws.Range("AB15").Value = 7.5
$ws.Range("AC15").Value = 13
$ws.Range("AD15").Value = 34
$ws.Range("AE15").Value = 11
$ws.Range("AF15").Value = 10
$ws.Range("AI15").Value = 19
$ws.Range("AN15").Value = 67
$ws.Range("AP15").Value = 1.44
$ws.Range("AQ15").Value = 2.7
$ws.Range("AR15").Value = 2.55
$ws.Range("AS15").Value = 1.49

# Row 17
$ws.Range("G17").Value = 2.5
$ws.Range("H17").Value = 3.1
$ws.Range("I17").Value = 2.35
$ws.Range("J17").Value = 3.2
$ws.Range("L17").Value = 3.1
$ws.Range("O17").Value = 1.36
$ws.Range("P17").Value = 3
$ws.Range("Q17").Value = 2.25
$ws.Range("R17").Value = 1.62
$ws.Range("S17").Value = 4
$ws.Range("T17").Value = 1.22
$ws.Range("Y17").Value = 8.5
$ws.Range("Z17").Value = 15
$ws.Range("AA17").Value = 12
$ws.Range("AB17").Value = 34
$ws.Range("AC17").Value = 29
$ws.Range("AE17").Value = 8
$ws.Range("AI17").Value = 7
$ws.Range("AJ17").Value = 11
$ws.Range("AK17").Value = 10
$ws.Range("AL17").Value = 21
$ws.Range("AM17").Value = 21
